# Sicom.xlsx amendment:
#  - "Combo" -> "COMBO" (uppercased) wherever it already appears (F2, F3),
#    and the same "COMBO" value added to the new F4/F5/F6 cells.
#  - "*Tndr Clsc, Spicy" -> "*Tndr Clsc, SPICY" (uppercased) at J4.
#  - Leaves the active selection on F6 (last cell touched by the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the existing "Combo" entries and extend the SelectMealOption
# column (F) with the same value for rows 4-6.
$ws.Range("F2").Value = "COMBO"
$ws.Range("F3").Value = "COMBO"
$ws.Range("F4").Value = "COMBO"
$ws.Range("F5").Value = "COMBO"
$ws.Range("F6").Value = "COMBO"

# Uppercase "Spicy" -> "SPICY" in the Flavor2 column for row 4.
$ws.Range("J4").Value = "*Tndr Clsc, SPICY"

# Move the active cell/selection to F6, matching the saved workbook state.
$ws.Range("F6").Select()
